# Automatische test-sync: 2025-06-17 22:49:11
# Appends a new "Afmelding nieuwsbrief" log entry to the Logs sheet,
# extends the conditional formatting ranges to cover the new row,
# and updates the Dashboard's "Afmelding" tally.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append new row 58 -----------------------------------
$ws = $wb.Worksheets.Item("Logs")

$newRow = 58
$ws.Cells.Item($newRow, 1).Value = "Afmelding nieuwsbrief"
$ws.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item($newRow, 3).Value = "Graag afmelden voor de nieuwsbrief. Dank u."
$ws.Cells.Item($newRow, 4).Value = "Afmelding"
$ws.Cells.Item($newRow, 6).Value = "2025-06-17 22:48:09"
$ws.Cells.Item($newRow, 7).Value = "Nee"

# --- Extend conditional formatting to include the new row ------------
$dFormats = $ws.Range("D2:D57").FormatConditions
for ($i = 1; $i -le $dFormats.Count; $i++) {
    $dFormats.Item($i).ModifyAppliesToRange($ws.Range("D2:D58"))
}

$gFormats = $ws.Range("G2:G57").FormatConditions
for ($i = 1; $i -le $gFormats.Count; $i++) {
    $gFormats.Item($i).ModifyAppliesToRange($ws.Range("G2:G58"))
}

# --- Dashboard sheet: bump the "Afmelding" count from 10 to 11 -------
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Cells.Item(4, 2).Value = 11
